$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header: "Tgl Masuk" (R1), bold like the other headers
$ws.Range("R1").Value = "Tgl Masuk"
$ws.Range("R1").Font.Bold = $true

# Re-affirm the "Kondisi" dropdown helper values in column G
$ws.Range("G2").Value = "NEW"
$ws.Range("G3").Value = "SECOND"

# New helper values for warranty day options (columns O and P)
$ws.Range("O2").Value = 365
$ws.Range("P2").Value = 365
$ws.Range("O3").Value = 30
$ws.Range("P3").Value = 30
$ws.Range("O4").Value = 7
$ws.Range("P4").Value = 7

# Sample "Tgl Masuk" (date entered) value, entered as a quoted/text date
$ws.Range("R2").Value = "'2023-12-30 00:00:00"
$ws.Range("R2").NumberFormat = "m/d/yy h:mm"

# Column width for the new column (closest achievable to 19.21875 chars
# given the engine's pixel-rounding on ColumnWidth)
$ws.Columns("R").ColumnWidth = 18.3

# Update selection to match the saved view state
$ws.Range("E10").Select() | Out-Null
